# Update "想去人数" (F column) figures across the four sheets to the
# newly scraped values, as per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 120
$ws.Cells.Item(3, 6).Value = 1283
$ws.Cells.Item(4, 6).Value = 923
$ws.Cells.Item(5, 6).Value = 969
$ws.Cells.Item(6, 6).Value = 1718
$ws.Cells.Item(8, 6).Value = 1145
$ws.Cells.Item(11, 6).Value = 109
$ws.Cells.Item(12, 6).Value = 261
$ws.Cells.Item(13, 6).Value = 36
$ws.Cells.Item(15, 6).Value = 641
$ws.Cells.Item(16, 6).Value = 132
$ws.Cells.Item(17, 6).Value = 89
$ws.Cells.Item(20, 6).Value = 320
$ws.Cells.Item(21, 6).Value = 103
$ws.Cells.Item(24, 6).Value = 626
$ws.Cells.Item(25, 6).Value = 133
$ws.Cells.Item(27, 6).Value = 838
$ws.Cells.Item(29, 6).Value = 112
$ws.Cells.Item(30, 6).Value = 24
$ws.Cells.Item(33, 6).Value = 9

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(4, 6).Value = 307
$ws.Cells.Item(7, 6).Value = 242
$ws.Cells.Item(11, 6).Value = 114

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 297

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 297
$ws.Cells.Item(3, 6).Value = 120
$ws.Cells.Item(4, 6).Value = 1284
$ws.Cells.Item(5, 6).Value = 923
$ws.Cells.Item(6, 6).Value = 969
$ws.Cells.Item(7, 6).Value = 1718
$ws.Cells.Item(9, 6).Value = 1145
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(13, 6).Value = 109
$ws.Cells.Item(14, 6).Value = 261
$ws.Cells.Item(15, 6).Value = 36
$ws.Cells.Item(17, 6).Value = 641
$ws.Cells.Item(18, 6).Value = 132
$ws.Cells.Item(19, 6).Value = 89
$ws.Cells.Item(22, 6).Value = 307
$ws.Cells.Item(25, 6).Value = 320
$ws.Cells.Item(27, 6).Value = 242
$ws.Cells.Item(28, 6).Value = 242
$ws.Cells.Item(29, 6).Value = 103
$ws.Cells.Item(32, 6).Value = 626
$ws.Cells.Item(33, 6).Value = 133
$ws.Cells.Item(35, 6).Value = 838
$ws.Cells.Item(39, 6).Value = 112
$ws.Cells.Item(40, 6).Value = 24
$ws.Cells.Item(43, 6).Value = 114
$ws.Cells.Item(44, 6).Value = 114
$ws.Cells.Item(46, 6).Value = 9
